$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:40:44.720041",
    "2021-10-05 13:40:44.720051",
    "2021-10-05 13:40:44.720055",
    "2021-10-05 13:40:44.720058",
    "2021-10-05 13:40:44.720061",
    "2021-10-05 13:40:44.720064",
    "2021-10-05 13:40:44.720066",
    "2021-10-05 13:40:44.720069",
    "2021-10-05 13:40:44.720072",
    "2021-10-05 13:40:44.720074",
    "2021-10-05 13:40:44.720077",
    "2021-10-05 13:40:44.720080",
    "2021-10-05 13:40:44.720083",
    "2021-10-05 13:40:44.720086",
    "2021-10-05 13:40:44.720088",
    "2021-10-05 13:40:44.720093",
    "2021-10-05 13:40:44.720098",
    "2021-10-05 13:40:44.720103",
    "2021-10-05 13:40:44.720108",
    "2021-10-05 13:40:44.720112"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
